$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16. This shifts the existing row 16
# (2a nueva(o) / Provincia de Melipilla) down to row 17, and the
# existing row 17 (1a nueva(o) / Peru) down to row 18.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the latest weekly price entry.
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 44476
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100112045
$ws.Cells.Item(16, 7).Value = "Zapallo"
$ws.Cells.Item(16, 8).Value = "Camote"
$ws.Cells.Item(16, 9).Value = "1a nueva(o)"
$ws.Cells.Item(16, 10).Value = 1200
$ws.Cells.Item(16, 11).Value = 480
$ws.Cells.Item(16, 12).Value = 500
$ws.Cells.Item(16, 13).Value = 490
$ws.Cells.Item(16, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 490
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = "Hortaliza"
